$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.282.28"
$ws.Range("E2").Value = "  -3.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.976.34"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.06"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.37"
$ws.Range("E7").Value = "  -10.60%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.27"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0828"
$ws.Range("E11").Value = "  +7.39%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.61"
$ws.Range("E13").Value = "  +4.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.861"
$ws.Range("E14").Value = "  -7.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.97"
$ws.Range("E15").Value = "  -6.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.267.36"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.981.87"
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.223.98"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.13"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0872"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.93"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  -4.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.76"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.79"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.132"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.89"
$ws.Range("E33").Value = "  -6.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0634"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -6.69%  "
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -8.99%  "
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.23"
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0966"
$ws.Range("E42").Value = "  -6.85%  "
$ws.Range("E43").Value = "  -4.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0213"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  -5.54%  "
$ws.Range("E46").Value = "  -8.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.21"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.47"
$ws.Range("E48").Value = "  -5.99%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.365.27"
$ws.Range("E49").Value = "  -3.51%  "
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("E51").Value = "  -3.13%  "
